# Update "想去人数" (F column) values across sheets to reflect the newly
# generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14993
$ws1.Range("F3").Value = 18951
$ws1.Range("F14").Value = 150
$ws1.Range("F17").Value = 1452
$ws1.Range("F20").Value = 96
$ws1.Range("F21").Value = 236
$ws1.Range("F22").Value = 7885
$ws1.Range("F24").Value = 32
$ws1.Range("F25").Value = 1
$ws1.Range("F27").Value = 1239
$ws1.Range("F29").Value = 6033
$ws1.Range("F31").Value = 70
$ws1.Range("F34").Value = 277
$ws1.Range("F35").Value = 5404
$ws1.Range("F36").Value = 63

# --- Sheet "演出" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 15

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14993
$ws4.Range("F3").Value = 18952
$ws4.Range("F14").Value = 150
$ws4.Range("F17").Value = 1452
$ws4.Range("F21").Value = 96
$ws4.Range("F22").Value = 236
$ws4.Range("F23").Value = 7885
$ws4.Range("F25").Value = 32
$ws4.Range("F26").Value = 1
$ws4.Range("F28").Value = 1239
$ws4.Range("F30").Value = 15
$ws4.Range("F32").Value = 6033
$ws4.Range("F34").Value = 70
$ws4.Range("F37").Value = 277
$ws4.Range("F38").Value = 5404
$ws4.Range("F39").Value = 64
